$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows for account 004515548 (FLAVIA) and 000772433 (MARCELO)
# These are data rows 239 and 240 (1-based, including the header row).
$ws.Range("A239:C240").EntireRow.Delete()
